# Auto-generated COM-interop script applying the "new word and rules 14 end season" commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a stray translation: G166 held the English placeholder "years" instead of
# the Russian translation "годы" (row is the adjective 작년 / "last year").
$ws.Cells.Item(166, 7).Value = "годы"

# --- Append the new "future plans" vocabulary block (rows 493-520) ---
# row 493
$ws.Cells.Item(493, 1).Value = "미래"
$ws.Cells.Item(493, 2).Value = "будущее"
$ws.Cells.Item(493, 4).Value = "time"
$ws.Cells.Item(493, 5).Value = "время"
$ws.Cells.Item(493, 6).Value = "relation"
$ws.Cells.Item(493, 7).Value = "относительно"
$ws.Cells.Item(493, 10).Value = "noun"
$ws.Cells.Item(493, 11).Value = 2
$ws.Cells.Item(493, 12).Value = 14
$ws.Cells.Item(493, 13).Value = "words"

# row 494
$ws.Cells.Item(494, 1).Value = "한국 진구를 라귀다"
$ws.Cells.Item(494, 2).Value = "встречаться с корейскими друзьями"
$ws.Cells.Item(494, 4).Value = "another"
$ws.Cells.Item(494, 5).Value = "другое"
$ws.Cells.Item(494, 11).Value = 2
$ws.Cells.Item(494, 12).Value = 14
$ws.Cells.Item(494, 13).Value = "expression"

# row 495
$ws.Cells.Item(495, 1).Value = "한국으로 여행을 가다"
$ws.Cells.Item(495, 2).Value = "отправляться в путешествие в Корею"
$ws.Cells.Item(495, 4).Value = "another"
$ws.Cells.Item(495, 5).Value = "другое"
$ws.Cells.Item(495, 11).Value = 2
$ws.Cells.Item(495, 12).Value = 14
$ws.Cells.Item(495, 13).Value = "expression"

# row 496
$ws.Cells.Item(496, 1).Value = "한국 회사에 쥐직하다"
$ws.Cells.Item(496, 2).Value = "устраиваться на работу в Корейскую компанию"
$ws.Cells.Item(496, 4).Value = "another"
$ws.Cells.Item(496, 5).Value = "другое"
$ws.Cells.Item(496, 11).Value = 2
$ws.Cells.Item(496, 12).Value = 14
$ws.Cells.Item(496, 13).Value = "expression"

# row 497
$ws.Cells.Item(497, 1).Value = "한국과 무역을 하다"
$ws.Cells.Item(497, 2).Value = "заниматься торговлей с Кореей"
$ws.Cells.Item(497, 4).Value = "another"
$ws.Cells.Item(497, 5).Value = "другое"
$ws.Cells.Item(497, 11).Value = 2
$ws.Cells.Item(497, 12).Value = 14
$ws.Cells.Item(497, 13).Value = "expression"

# row 498
$ws.Cells.Item(498, 1).Value = "한국 문화에 관심이 있다"
$ws.Cells.Item(498, 2).Value = "интересоваться Корейской культурой"
$ws.Cells.Item(498, 4).Value = "another"
$ws.Cells.Item(498, 5).Value = "другое"
$ws.Cells.Item(498, 11).Value = 2
$ws.Cells.Item(498, 12).Value = 14
$ws.Cells.Item(498, 13).Value = "expression"

# row 499
$ws.Cells.Item(499, 1).Value = "한국으로 유학을 가다"
$ws.Cells.Item(499, 2).Value = "отправляться на стажировуц в Корею"
$ws.Cells.Item(499, 4).Value = "another"
$ws.Cells.Item(499, 5).Value = "другое"
$ws.Cells.Item(499, 11).Value = 2
$ws.Cells.Item(499, 12).Value = 14
$ws.Cells.Item(499, 13).Value = "expression"

# row 500
$ws.Cells.Item(500, 1).Value = "번역가가 되다"
$ws.Cells.Item(500, 2).Value = "становиться письменным переводчиком"
$ws.Cells.Item(500, 4).Value = "another"
$ws.Cells.Item(500, 5).Value = "другое"
$ws.Cells.Item(500, 11).Value = 2
$ws.Cells.Item(500, 12).Value = 14
$ws.Cells.Item(500, 13).Value = "expression"

# row 501
$ws.Cells.Item(501, 1).Value = "통역사가 되다"
$ws.Cells.Item(501, 2).Value = "становиться устным переводчиком"
$ws.Cells.Item(501, 4).Value = "another"
$ws.Cells.Item(501, 5).Value = "другое"
$ws.Cells.Item(501, 11).Value = 2
$ws.Cells.Item(501, 12).Value = 14
$ws.Cells.Item(501, 13).Value = "expression"

# row 502
$ws.Cells.Item(502, 1).Value = "선생님이 되다"
$ws.Cells.Item(502, 2).Value = "становиться учителем"
$ws.Cells.Item(502, 4).Value = "another"
$ws.Cells.Item(502, 5).Value = "другое"
$ws.Cells.Item(502, 11).Value = 2
$ws.Cells.Item(502, 12).Value = 14
$ws.Cells.Item(502, 13).Value = "expression"

# row 503
$ws.Cells.Item(503, 1).Value = "꿈을 이루다"
$ws.Cells.Item(503, 2).Value = "достигать мечту"
$ws.Cells.Item(503, 4).Value = "another"
$ws.Cells.Item(503, 5).Value = "другое"
$ws.Cells.Item(503, 11).Value = 2
$ws.Cells.Item(503, 12).Value = 14
$ws.Cells.Item(503, 13).Value = "expression"

# row 504
$ws.Cells.Item(504, 1).Value = "좋은 사람과 결흔하다"
$ws.Cells.Item(504, 2).Value = "играть свадьбу с хорошим человеком"
$ws.Cells.Item(504, 4).Value = "another"
$ws.Cells.Item(504, 5).Value = "другое"
$ws.Cells.Item(504, 11).Value = 2
$ws.Cells.Item(504, 12).Value = 14
$ws.Cells.Item(504, 13).Value = "expression"

# row 505
$ws.Cells.Item(505, 1).Value = "성공하다"
$ws.Cells.Item(505, 2).Value = "добиваться успеха"
$ws.Cells.Item(505, 4).Value = "another"
$ws.Cells.Item(505, 5).Value = "другое"
$ws.Cells.Item(505, 11).Value = 2
$ws.Cells.Item(505, 12).Value = 14
$ws.Cells.Item(505, 13).Value = "expression"

# row 506
$ws.Cells.Item(506, 1).Value = "돈을 많이 벌다"
$ws.Cells.Item(506, 2).Value = "зарабатывать много денег"
$ws.Cells.Item(506, 4).Value = "another"
$ws.Cells.Item(506, 5).Value = "другое"
$ws.Cells.Item(506, 11).Value = 2
$ws.Cells.Item(506, 12).Value = 14
$ws.Cells.Item(506, 13).Value = "expression"

# row 507
$ws.Cells.Item(507, 1).Value = "외국 여행을 가다"
$ws.Cells.Item(507, 2).Value = "отправляться в путешествие за границу"
$ws.Cells.Item(507, 4).Value = "another"
$ws.Cells.Item(507, 5).Value = "другое"
$ws.Cells.Item(507, 11).Value = 2
$ws.Cells.Item(507, 12).Value = 14
$ws.Cells.Item(507, 13).Value = "expression"

# row 508
$ws.Cells.Item(508, 1).Value = "무역"
$ws.Cells.Item(508, 2).Value = "торговля"
$ws.Cells.Item(508, 4).Value = "action"
$ws.Cells.Item(508, 5).Value = "действие"
$ws.Cells.Item(508, 10).Value = "noun"
$ws.Cells.Item(508, 11).Value = 2
$ws.Cells.Item(508, 12).Value = 14
$ws.Cells.Item(508, 13).Value = "words"

# row 509
$ws.Cells.Item(509, 1).Value = "양복"
$ws.Cells.Item(509, 2).Value = "деловой костюм"
$ws.Cells.Item(509, 4).Value = "clothes"
$ws.Cells.Item(509, 5).Value = "одежда"
$ws.Cells.Item(509, 10).Value = "noun"
$ws.Cells.Item(509, 11).Value = 2
$ws.Cells.Item(509, 12).Value = 14
$ws.Cells.Item(509, 13).Value = "words"

# row 510
$ws.Cells.Item(510, 1).Value = "사업을 하다"
$ws.Cells.Item(510, 2).Value = "заниматься предпринимательством"
$ws.Cells.Item(510, 4).Value = "action"
$ws.Cells.Item(510, 5).Value = "действие"
$ws.Cells.Item(510, 11).Value = 2
$ws.Cells.Item(510, 12).Value = 14
$ws.Cells.Item(510, 13).Value = "expression"

# row 511
$ws.Cells.Item(511, 1).Value = "한 살"
$ws.Cells.Item(511, 2).Value = "один год"
$ws.Cells.Item(511, 4).Value = "time"
$ws.Cells.Item(511, 5).Value = "время"
$ws.Cells.Item(511, 6).Value = "years"
$ws.Cells.Item(511, 7).Value = "год"
$ws.Cells.Item(511, 11).Value = 2
$ws.Cells.Item(511, 12).Value = 14
$ws.Cells.Item(511, 13).Value = "expression"

# row 512
$ws.Cells.Item(512, 1).Value = "신청하다"
$ws.Cells.Item(512, 2).Value = "подавать заявку"
$ws.Cells.Item(512, 4).Value = "action"
$ws.Cells.Item(512, 5).Value = "действие"
$ws.Cells.Item(512, 11).Value = 2
$ws.Cells.Item(512, 12).Value = 14
$ws.Cells.Item(512, 13).Value = "expression"

# row 513
$ws.Cells.Item(513, 1).Value = "가사"
$ws.Cells.Item(513, 2).Value = "текст"
$ws.Cells.Item(513, 4).Value = "another"
$ws.Cells.Item(513, 5).Value = "другое"
$ws.Cells.Item(513, 10).Value = "noun"
$ws.Cells.Item(513, 11).Value = 2
$ws.Cells.Item(513, 12).Value = 14
$ws.Cells.Item(513, 13).Value = "words"

# row 514
$ws.Cells.Item(514, 1).Value = "앉다"
$ws.Cells.Item(514, 2).Value = "садиться"
$ws.Cells.Item(514, 4).Value = "action"
$ws.Cells.Item(514, 5).Value = "действие"
$ws.Cells.Item(514, 10).Value = "verb"
$ws.Cells.Item(514, 11).Value = 2
$ws.Cells.Item(514, 12).Value = 14
$ws.Cells.Item(514, 13).Value = "words"

# row 515
$ws.Cells.Item(515, 1).Value = "아이들"
$ws.Cells.Item(515, 2).Value = "дети"
$ws.Cells.Item(515, 4).Value = "people"
$ws.Cells.Item(515, 5).Value = "люди"
$ws.Cells.Item(515, 10).Value = "noun"
$ws.Cells.Item(515, 11).Value = 2
$ws.Cells.Item(515, 12).Value = 14
$ws.Cells.Item(515, 13).Value = "words"

# row 516
$ws.Cells.Item(516, 1).Value = "부자"
$ws.Cells.Item(516, 2).Value = "богатый человек"
$ws.Cells.Item(516, 4).Value = "people"
$ws.Cells.Item(516, 5).Value = "люди"
$ws.Cells.Item(516, 10).Value = "noun"
$ws.Cells.Item(516, 11).Value = 2
$ws.Cells.Item(516, 12).Value = 14
$ws.Cells.Item(516, 13).Value = "words"

# row 517
$ws.Cells.Item(517, 1).Value = "전 세계"
$ws.Cells.Item(517, 2).Value = "весь мир"
$ws.Cells.Item(517, 4).Value = "location"
$ws.Cells.Item(517, 5).Value = "место"
$ws.Cells.Item(517, 11).Value = 2
$ws.Cells.Item(517, 12).Value = 14
$ws.Cells.Item(517, 13).Value = "expression"

# row 518
$ws.Cells.Item(518, 1).Value = "초급"
$ws.Cells.Item(518, 2).Value = "начальный уровень"
$ws.Cells.Item(518, 4).Value = "another"
$ws.Cells.Item(518, 5).Value = "другое"
$ws.Cells.Item(518, 11).Value = 2
$ws.Cells.Item(518, 12).Value = 14
$ws.Cells.Item(518, 13).Value = "expression"

# row 519
$ws.Cells.Item(519, 1).Value = "중급"
$ws.Cells.Item(519, 2).Value = "средний уровень"
$ws.Cells.Item(519, 4).Value = "another"
$ws.Cells.Item(519, 5).Value = "другое"
$ws.Cells.Item(519, 11).Value = 2
$ws.Cells.Item(519, 12).Value = 14
$ws.Cells.Item(519, 13).Value = "expression"

# row 520
$ws.Cells.Item(520, 1).Value = "고급"
$ws.Cells.Item(520, 2).Value = "высший уровень"
$ws.Cells.Item(520, 4).Value = "another"
$ws.Cells.Item(520, 5).Value = "другое"
$ws.Cells.Item(520, 11).Value = 2
$ws.Cells.Item(520, 12).Value = 14
$ws.Cells.Item(520, 13).Value = "expression"

# --- Column B was manually narrowed (author dragged the border), which also clears the
# "best fit" auto-size flag. 26.5 is the closest value this host's ColumnWidth setter
# (quantized to 1/6-character steps) can reach to the authored 27.33203125 raw width.
$ws.Columns.Item(2).ColumnWidth = 26.5

# --- Restore the sheet selection to the cell the author left active (F166) ---
$ws.Range("F166").Select()
